# Module 9 assignment: sort the incarceration-rate data by
# "Population_per_100,000" (column E) in descending order, then leave the
# selection on the last data row of that column (E16) as Excel does after
# running the Data > Sort command from a selection inside that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Whole data block (including header row) and the sort key column.
$sortRange = $ws.Range("A1:H16")
$keyRange  = $ws.Range("E1:E16")

# Key1, Order1=xlDescending(2), ..., Header=xlYes(1)
$sortRange.Sort($keyRange, 2, $null, $null, $null, $null, $null, 1)

# Match the final selection left behind in the saved workbook.
[void]$ws.Range("E16").Select()
